# ---------------------------------------------------------------------------
# Design updates and images
#   - refresh the cached "datetimeFigureOut" footer field (slide master +
#     every slide layout) from 31/03/2017 to 01/05/2017
#   - reposition / re-center the three text boxes on the single content
#     slide, fill in the "Current Score" line, and collapse the
#     "? ? ? ? Put" run into a single "_ _ _ _     Put" run
# ---------------------------------------------------------------------------

# Shape.Left/.Top are exposed in points but the engine rounds through a
# single-precision float before converting back to EMU, which can lose the
# last EMU of a large offset. Search for the nearest point value whose
# round-trip lands exactly on the EMU we want so the emitted <a:off> matches
# the source OOXML exactly.
function EmuToPt($emu) {
    $base = $emu / 12700.0
    for ($i = 0; $i -lt 500; $i++) {
        $cand = $base + ($i * 0.0000001)
        $f32 = [single]$cand
        $backEmu = [math]::Floor([double]$f32 * 12700.0)
        if ($backEmu -eq $emu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Date placeholder text (slide master + all custom layouts)
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "01/05/2017"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 1 shape tweaks
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# --- "TextBox 1" (Topic / Total words / Current word / Current Score) -----
$topicBox = $s.Shapes.Item("TextBox 1")
$topicBox.Left = EmuToPt 3837393
$topicBox.Top = EmuToPt 1050527

$topicRange = $topicBox.TextFrame.TextRange

# Add the "Current Score: 25" text as a new paragraph that inherits the
# formatting of the preceding "Current word: 1" paragraph, then drop the
# old trailing empty paragraph that got pushed down by the insertion.
$thirdPara = $topicRange.Paragraphs(3, 1)
$cr = [char]13
$thirdPara.InsertAfter($cr + "Current Score: 25") | Out-Null
$topicRange.Paragraphs(5, 1).Delete() | Out-Null

# Center every paragraph in the box.
for ($i = 1; $i -le 4; $i++) {
    $topicRange.Paragraphs($i, 1).ParagraphFormat.Alignment = 2
}

# --- "Rectangle 2" (answer bar background) ---------------------------------
$rect2 = $s.Shapes.Item("Rectangle 2")
$rect2.Left = EmuToPt 2975294
$rect2.Top = EmuToPt 4772763

# --- "TextBox 14" ("? ? ? ?  Put" -> "_ _ _ _     Put") --------------------
$putBox = $s.Shapes.Item("TextBox 14")
$putBox.Left = EmuToPt 3400425
$putBox.Top = EmuToPt 3568526

$putRange = $putBox.TextFrame.TextRange
$newText = "_ _ _ _     Put"

$firstRun = $putRange.Runs(1, 1)
$firstRun.Text = $newText

$tailStart = $firstRun.Start + $newText.Length
$tailLength = $putRange.Text.Length - ($tailStart - 1)
if ($tailLength -gt 0) {
    $putRange.Characters($tailStart, $tailLength).Delete() | Out-Null
}

$putRange.Paragraphs(1, 1).ParagraphFormat.Alignment = 2
